# Apply the "Ajout de l'authorPerson" change to StructureDefinition-Author.xlsx
$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# 1. Update the Date value (Metadata!B8)
$wsMetadata.Range("B8").Value = "2025-05-02T17:24:18+00:00"

# 2. Append a sentence to the Description text
#    (Metadata!B12, and Elements!M2 which shares the same string)
$newDescription = "Modèle logique d'un auteur.`nReprésente les personnes physiques et/ou les systèmes (dispositifs, automates, services numériques référencés…) auteurs d’un document ou d'un lot de soummission. `nL’auteur peut être : `n- Un professionnel (personne physique) via son logiciel de professionnel, `n- Le patient/usager (personne physique) via Mon espace Santé, `n- Un système de structure (dispositif, automate, appareil connecté…), `n- Un SNR (Service Numérique Référencé), `nauthor est un ensemble constitué des sous-attributs authorInstitution, authorPerson, authorRole et authorSpecialty et ne porte pas de valeur par lui-même. "
$wsMetadata.Range("B12").Value = $newDescription
$wsElements.Range("M2").Value = $newDescription

# 3. Replace the Author.person Type(s) value with the new Reference(...) text (Elements!K4)
$newAuthorPersonType = "Reference(https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/AuthorPersonPS|https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/AuthorPersonPatient|https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/AuthorPersonSNR|https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/AuthorPersonSystem)`n"
$wsElements.Range("K4").Value = $newAuthorPersonType

# 4. Widen column K (11) on the Elements sheet to fit the new, longer content
$wsElements.Columns.Item(11).ColumnWidth = 254.16666666666666
